$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.819.78"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +4.81%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.113.82"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.07%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.50"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.91%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.16"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.58%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.107.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.26%  "

$ws.Range("E9").Value = "  +1.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.149"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +10.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.76"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +6.92%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.467"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.25%  "

$ws.Range("E13").Value = "  +6.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.58"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.19%  "

$ws.Range("E15").Value = "  -0.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.629.71"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +3.13%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.16"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.50%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.111.47"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "62.737.44"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.78%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "468.22"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +6.29%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.07"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.729"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.56"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +5.89%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.34"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.17"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.24"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.68"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.28"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +4.75%  "

$ws.Range("E30").Value = "  -0.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.80"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +7.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.02"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.55%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.110"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +3.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0875"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +9.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.38"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +11.98%  "

$ws.Range("E36").Value = "  +3.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.05"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.29"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +16.34%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.94"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "435.09"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +6.92%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.75"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.932.13"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.76%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0370"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.278"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +9.31%  "

$ws.Range("E45").Value = "  +3.16%  "

$ws.Range("E46").Value = "  +5.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.68"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +4.47%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.70"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.23%  "

$ws.Range("E50").Value = "  +0.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.63"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.49%  "
